$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update STOCK_QUANTITY values to reconnect backend data with frontend
$ws.Range("E2").Value = 100
$ws.Range("E3").Value = 77

# Move the active selection/cursor to E5
$ws.Range("E5").Select()
